$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("edges")

# Update the value in E2 from 0.2 to 0.5
$ws.Range("E2").Value = 0.5

# Update the selected cell/range on the sheet from E3 to E2
$ws.Activate()
$ws.Range("E2").Select()
